$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.938.75"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.817.41"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.75"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4292"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3695"
$ws.Range("E8").Value = "  +2.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07245"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8638"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").Value = "2.026.55"
$ws.Range("E11").Value = "  +11.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.06"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.649"
$ws.Range("E13").Value = "  +3.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.386"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06924"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.69"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008919"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.20"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "26.991.07"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.186"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "2.262.98"
$ws.Range("E24").Value = "  +11.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.95"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.887"
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.24"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.212"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.890"
$ws.Range("E29").Value = "  +16.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.01"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08960"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.412"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.804"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05213"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01925"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5075"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.738"
$ws.Range("E41").Value = "  +8.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1641"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.453"
$ws.Range("E43").Value = "  +7.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.225"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.07"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.42"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.656"
$ws.Range("E48").Value = "  +4.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06305"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4561"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.804"
$ws.Range("E51").Value = "  +4.67%  "

# Row 32 becomes ARBITRUM data, Row 33 becomes ImmutableX data
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.167"
$ws.Range("E32").Value = "  +6.48%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7429"
$ws.Range("E33").Value = "  +2.25%  "
